# Actualización automática 2025-11-27 17:30:08
$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": update product-group breakdown for
#     client "AGUILAR REYES CESAR VINICIO" (row 4) for the new sale data
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("D4").Value = 907.6799999999999
$wsGrupo.Range("E4").Value = 139.11
$wsGrupo.Range("M4").Value = 9443.559999999999

# --- Sheet "VENTA MENSUAL": update the november (noviembre) sale total
#     for the same client (row 4); the grand total row (60) holds the
#     column total as a static value, so it must be updated too
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F4").Value = 12411.83
$wsMensual.Range("F60").Value = 88209.84

# --- Sheet "CUMPLIMIENTO MENSUAL": update VENTA (D), POR CUMPLIR (E)
#     and CUMPLIMIENTO (F) for the affected groups, and the TOTAL row
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3: 240X80 PORCELANATO
$wsCumplimiento.Range("D3").Value = 10967.52
$wsCumplimiento.Range("E3").Value = -4344.26
$wsCumplimiento.Range("F3").Value = 1.655909627585207

# Row 4: FREGADEROS DE COCINA
$wsCumplimiento.Range("D4").Value = 352.48
$wsCumplimiento.Range("E4").Value = 492.22
$wsCumplimiento.Range("F4").Value = 0.4172842429264828

# Row 12: PORCELANATO
$wsCumplimiento.Range("D12").Value = 45468.3
$wsCumplimiento.Range("E12").Value = 19475.7
$wsCumplimiento.Range("F12").Value = 0.7001154841093866

# Row 14: TOTAL
$wsCumplimiento.Range("D14").Value = 84964.59
$wsCumplimiento.Range("E14").Value = 13991.66685923838
$wsCumplimiento.Range("F14").Value = 0.858607557487335

$wb.Save()
